$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("L2").Value = 1.4
$ws.Range("M2").Value = 3
$ws.Range("U2").Value = 10
# Row 3
$ws.Range("G3").Value = 1.67
$ws.Range("H3").Value = 3.4
$ws.Range("I3").Value = 5.75
$ws.Range("J3").Value = 1.08
$ws.Range("K3").Value = 8
$ws.Range("L3").Value = 1.44
$ws.Range("M3").Value = 2.63
$ws.Range("N3").Value = 2.4
$ws.Range("O3").Value = 1.53
$ws.Range("R3").Value = 2.5
$ws.Range("S3").Value = 1.5
$ws.Range("U3").Value = 6.5
$ws.Range("W3").Value = 12
$ws.Range("X3").Value = 17
$ws.Range("AA3").Value = 7
$ws.Range("AE3").Value = 11
$ws.Range("AF3").Value = 26
$ws.Range("AG3").Value = 21
$ws.Range("AH3").Value = 67
$ws.Range("AJ3").Value = 67
# Row 4
$ws.Range("G4").Value = 2.38
$ws.Range("P4").Value = 1.57
# Row 5
$ws.Range("G5").Value = 2.1
$ws.Range("H5").Value = 3.2
$ws.Range("I5").Value = 3.7
$ws.Range("K5").Value = 7.5
$ws.Range("N5").Value = 2.25
$ws.Range("O5").Value = 1.62
$ws.Range("Q5").Value = 2.5
$ws.Range("U5").Value = 9
$ws.Range("W5").Value = 19
$ws.Range("X5").Value = 19
$ws.Range("AE5").Value = 9
$ws.Range("AF5").Value = 17
$ws.Range("AH5").Value = 41
$ws.Range("AI5").Value = 34
# Row 6
$ws.Range("G6").Value = 1.7
$ws.Range("J6").Value = 1.06
$ws.Range("K6").Value = 10
$ws.Range("P6").Value = 1.4
$ws.Range("Q6").Value = 2.75
$ws.Range("R6").Value = 1.83
$ws.Range("S6").Value = 1.83
$ws.Range("T6").Value = 7
$ws.Range("Z6").Value = 10
$ws.Range("AG6").Value = 17
# Row 7
$ws.Range("G7").Value = 2.6
$ws.Range("I7").Value = 2.88
$ws.Range("J7").Value = 1.11
$ws.Range("K7").Value = 6.5
$ws.Range("P7").Value = 1.57
# Row 8
$ws.Range("G8").Value = 1.67
$ws.Range("P8").Value = 1.44
$ws.Range("Q8").Value = 2.63
$ws.Range("R8").Value = 2.1
$ws.Range("S8").Value = 1.67
$ws.Range("AH8").Value = 67
$ws.Range("AI8").Value = 51
# Row 11
$ws.Range("J11").Value = 1.08
$ws.Range("K11").Value = 8
$ws.Range("N11").Value = 2.3
$ws.Range("O11").Value = 1.6
# Row 12
$ws.Range("G12").Value = 1.57
$ws.Range("H12").Value = 3.55
$ws.Range("I12").Value = 5.8
$ws.Range("Q12").Value = 2.35
$ws.Range("R12").Value = 2.1
$ws.Range("T12").Value = 5.4
$ws.Range("U12").Value = 6.4
$ws.Range("W12").Value = 10.75
$ws.Range("X12").Value = 14.5
$ws.Range("Z12").Value = 7.9
$ws.Range("AA12").Value = 7.2
$ws.Range("AE12").Value = 12
$ws.Range("AF12").Value = 32
$ws.Range("AG12").Value = 20
$ws.Range("AH12").Value = 120
$ws.Range("AI12").Value = 80
$ws.Range("AJ12").Value = 90
# Row 13
$ws.Range("K13").Value = 10
# Row 19
$ws.Range("G19").Value = 1.85
$ws.Range("I19").Value = 3.75
$ws.Range("U19").Value = 9.5
$ws.Range("AE19").Value = 13
$ws.Range("AF19").Value = 21
# Row 20
$ws.Range("G20").Value = 1.2
$ws.Range("H20").Value = 7
$ws.Range("I20").Value = 10
$ws.Range("K20").Value = 34
$ws.Range("L20").Value = 1.05
$ws.Range("M20").Value = 8.5
$ws.Range("N20").Value = 1.2
$ws.Range("O20").Value = 4.33
$ws.Range("R20").Value = 1.57
$ws.Range("S20").Value = 2.25
$ws.Range("T20").Value = 17
$ws.Range("U20").Value = 11
$ws.Range("V20").Value = 12
$ws.Range("W20").Value = 10
$ws.Range("AB20").Value = 19
$ws.Range("AD20").Value = 101
$ws.Range("AG20").Value = 29
$ws.Range("AJ20").Value = 41
# Row 21
$ws.Range("G21").Value = 1.62
$ws.Range("H21").Value = 4
$ws.Range("I21").Value = 4.75
$ws.Range("N21").Value = 1.5
$ws.Range("O21").Value = 2.5
$ws.Range("V21").Value = 9
$ws.Range("AA21").Value = 8
# Row 22
$ws.Range("G22").Value = 3.75
$ws.Range("J22").Value = 21
$ws.Range("K22").Value = 1.02
$ws.Range("Z22").Value = 21
$ws.Range("AA22").Value = 8.5
$ws.Range("AG22").Value = 9
# Row 23
$ws.Range("G23").Value = 2.12
$ws.Range("H23").Value = 3
$ws.Range("I23").Value = 3.5
$ws.Range("M23").Value = 2.27
$ws.Range("N23").Value = 2.4
$ws.Range("O23").Value = 1.44
$ws.Range("P23").Value = 1.52
$ws.Range("Q23").Value = 2.22
$ws.Range("R23").Value = 2.07
$ws.Range("T23").Value = 5.5
$ws.Range("U23").Value = 8.75
$ws.Range("V23").Value = 9.5
$ws.Range("W23").Value = 20
$ws.Range("X23").Value = 22
$ws.Range("Z23").Value = 6.4
$ws.Range("AA23").Value = 6
$ws.Range("AB23").Value = 19
$ws.Range("AE23").Value = 7.8
$ws.Range("AF23").Value = 17
$ws.Range("AG23").Value = 13
$ws.Range("AH23").Value = 50
$ws.Range("AI23").Value = 40
$ws.Range("AJ23").Value = 60
